$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 4)
$ws.Range("A4").Value = 44902
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B4").Value = "More prepping, getting towards temp ACF"
$ws.Range("C4").Value = 2

# Widen column B to fit the new, longer text
$ws.Range("B1").ColumnWidth = 45.65

# Update the active selection
$ws.Range("B12").Select()
